$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.Value = "'" + $value
    $c.Style = "Normal"
}

Set-TextCell 'D2' '27.133.28'
Set-TextCell 'E2' '  +1.08%  '
Set-TextCell 'D3' '1.891.15'
Set-TextCell 'E3' '  +1.90%  '
Set-TextCell 'D4' '0.9998'
Set-TextCell 'E4' '  -0.04%  '
Set-TextCell 'D5' '308.13'
Set-TextCell 'E5' '  +1.26%  '
Set-TextCell 'D6' '0.9995'
Set-TextCell 'E6' '  -0.09%  '
Set-TextCell 'D7' '0.5167'
Set-TextCell 'E7' '  +2.47%  '
Set-TextCell 'D8' '0.3720'
Set-TextCell 'E8' '  +1.91%  '
Set-TextCell 'D9' '0.07213'
Set-TextCell 'E9' '  +0.60%  '
Set-TextCell 'D10' '0.9044'
Set-TextCell 'E10' '  +1.53%  '
Set-TextCell 'D11' '21.04'
Set-TextCell 'E11' '  +2.03%  '
Set-TextCell 'D12' '0.07622'
Set-TextCell 'E12' '  +1.41%  '
Set-TextCell 'D13' '1.896.65'
Set-TextCell 'E13' '  +2.13%  '
Set-TextCell 'D14' '94.97'
Set-TextCell 'E14' '  +3.20%  '
Set-TextCell 'D15' '5.276'
Set-TextCell 'E15' '  +0.96%  '
Set-TextCell 'E16' '  -0.04%  '
Set-TextCell 'D17' '0.000008510'
Set-TextCell 'E17' '  +0.16%  '
Set-TextCell 'D18' '14.36'
Set-TextCell 'E18' '  +2.22%  '
Set-TextCell 'D19' '0.9994'
Set-TextCell 'E19' '  -0.09%  '
Set-TextCell 'D20' '27.166.79'
Set-TextCell 'D21' '5.056'
Set-TextCell 'E21' '  +0.61%  '
Set-TextCell 'D22' '2.131.81'
Set-TextCell 'E22' '  +2.04%  '
Set-TextCell 'D23' '10.58'
Set-TextCell 'E23' '  +2.51%  '
Set-TextCell 'D24' '6.436'
Set-TextCell 'E24' '  -0.27%  '
Set-TextCell 'D25' '145.15'
Set-TextCell 'E25' '  -0.89%  '
Set-TextCell 'D26' '1.789'
Set-TextCell 'E26' '  -0.54%  '
Set-TextCell 'D27' '18.08'
Set-TextCell 'E27' '  +1.45%  '
Set-TextCell 'E28' '  +4.85%  '
Set-TextCell 'D29' '114.68'
Set-TextCell 'E29' '  +1.59%  '
Set-TextCell 'D30' '4.983'
Set-TextCell 'E30' '  +7.01%  '
Set-TextCell 'D31' '4.815'
Set-TextCell 'E31' '  +3.89%  '
Set-TextCell 'D32' '0.09214'
Set-TextCell 'E32' '  +0.19%  '
Set-TextCell 'D33' '0.05056'
Set-TextCell 'E33' '  -0.62%  '
Set-TextCell 'D34' '1.198'
Set-TextCell 'E34' '  +4.62%  '
Set-TextCell 'D35' '0.7591'
Set-TextCell 'E35' '  +3.17%  '
Set-TextCell 'D36' '3.030'
Set-TextCell 'E36' '  +1.30%  '
Set-TextCell 'D37' '3.275'
Set-TextCell 'E37' '  +1.19%  '
Set-TextCell 'D38' '2.565'
Set-TextCell 'E38' '  +2.07%  '
Set-TextCell 'D39' '0.5652'
Set-TextCell 'E39' '  +6.10%  '
Set-TextCell 'D40' '0.01997'
Set-TextCell 'E40' '  +0.24%  '
Set-TextCell 'D41' '1.078'
Set-TextCell 'E41' '  +0.46%  '
Set-TextCell 'D42' '8.942'
Set-TextCell 'E42' '  +7.14%  '
Set-TextCell 'B43' 'FraxShare'
Set-TextCell 'C43' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell 'D43' '6.595'
Set-TextCell 'E43' '  +1.84%  '
Set-TextCell 'B44' 'Quant'
Set-TextCell 'C44' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell 'D44' '118.43'
Set-TextCell 'E44' '  -0.60%  '
Set-TextCell 'D45' '0.1510'
Set-TextCell 'E45' '  +3.01%  '
Set-TextCell 'D46' '0.4816'
Set-TextCell 'E46' '  +3.90%  '
Set-TextCell 'D47' '10.22'
Set-TextCell 'E47' '  +2.93%  '
Set-TextCell 'D48' '0.9994'
Set-TextCell 'E48' '  -0.09%  '
Set-TextCell 'D50' '37.18'
Set-TextCell 'E50' '  +0.82%  '
Set-TextCell 'D51' '63.58'
Set-TextCell 'E51' '  +1.24%  '
